$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 4638263.5
$ws.Cells.Item(18, 9).Value = 9259859
$ws.Cells.Item(18, 10).Value = 16668
$ws.Cells.Item(18, 11).Value = 9259859
$ws.Cells.Item(18, 12).Value = 16668
$ws.Cells.Item(18, 13).Value = -9259575
$ws.Cells.Item(18, 14).Value = -17236
# Row 113
$ws.Cells.Item(113, 8).Value = 5193.6665
$ws.Cells.Item(113, 9).Value = 4318.091
$ws.Cells.Item(113, 10).Value = 7601.5
$ws.Cells.Item(113, 11).Value = 4318.091
$ws.Cells.Item(113, 12).Value = 7601.5
$ws.Cells.Item(113, 13).Value = -1064.091
$ws.Cells.Item(113, 14).Value = -14109.5
# Row 116
$ws.Cells.Item(116, 8).Value = 2756.149
$ws.Cells.Item(116, 9).Value = 2216.1562
$ws.Cells.Item(116, 10).Value = 3908.1333
$ws.Cells.Item(116, 11).Value = 2216.1562
$ws.Cells.Item(116, 12).Value = 3908.1333
$ws.Cells.Item(116, 13).Value = 1225.8438
$ws.Cells.Item(116, 14).Value = -10792.1333
# Row 132
$ws.Cells.Item(132, 8).Value = 10533021
$ws.Cells.Item(132, 9).Value = 12507212
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 37521636
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).Value = -37519106
$ws.Cells.Item(132, 14).Value = -17060

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 13160403
$ws.Cells.Item(2, 9).Value = 16667843
$ws.Cells.Item(2, 10).Value = 7500
$ws.Cells.Item(2, 11).Value = 16667843
$ws.Cells.Item(2, 12).Value = 7500
$ws.Cells.Item(2, 13).Value = -16667730
$ws.Cells.Item(2, 14).Value = -7726
# Row 45
$ws.Cells.Item(45, 8).Value = 1581.125
$ws.Cells.Item(45, 9).Value = 1047.931
$ws.Cells.Item(45, 10).Value = 2986.818
$ws.Cells.Item(45, 11).Value = 1047.931
$ws.Cells.Item(45, 12).Value = 2986.818
$ws.Cells.Item(45, 13).Value = -670.931
$ws.Cells.Item(45, 14).Value = -3740.818
# Row 116
$ws.Cells.Item(116, 8).Value = 13160403
$ws.Cells.Item(116, 9).Value = 16667843
$ws.Cells.Item(116, 10).Value = 7500
$ws.Cells.Item(116, 11).Value = 16667843
$ws.Cells.Item(116, 12).Value = 7500
$ws.Cells.Item(116, 13).Value = -16665549
$ws.Cells.Item(116, 14).Value = -12088
# Row 132
$ws.Cells.Item(132, 8).Value = 27782126
$ws.Cells.Item(132, 9).Value = 45458840
$ws.Cells.Item(132, 10).Value = 4430.4287
$ws.Cells.Item(132, 11).Value = 136376520
$ws.Cells.Item(132, 12).Value = 13291.2861
$ws.Cells.Item(132, 13).Value = -136373990
$ws.Cells.Item(132, 14).Value = -18351.2861

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 13160403
$ws.Cells.Item(3, 9).Value = 16667843
$ws.Cells.Item(3, 10).Value = 7500
$ws.Cells.Item(3, 11).Value = 16667843
$ws.Cells.Item(3, 12).Value = 7500
$ws.Cells.Item(3, 13).Value = -16667729
$ws.Cells.Item(3, 14).Value = -7728
# Row 107
$ws.Cells.Item(107, 8).Value = 4333.3335
$ws.Cells.Item(107, 9).Value = 1500
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 13).Value = 420

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 144.61111
$ws.Cells.Item(7, 9).Value = 92.5
$ws.Cells.Item(7, 10).Value = 186.3
$ws.Cells.Item(7, 11).Value = 92.5
$ws.Cells.Item(7, 12).Value = 186.3
$ws.Cells.Item(7, 13).Value = 20.5
$ws.Cells.Item(7, 14).Value = -412.3
# Row 22
$ws.Cells.Item(22, 8).Value = 1012.1053
$ws.Cells.Item(22, 10).Value = 2232.8572
$ws.Cells.Item(22, 12).Value = 2232.8572
$ws.Cells.Item(22, 14).Value = -2932.8572
# Row 31
$ws.Cells.Item(31, 8).Value = 2634348.5
$ws.Cells.Item(31, 9).Value = 3227582.2
$ws.Cells.Item(31, 10).Value = 7171.143
$ws.Cells.Item(31, 11).Value = 3227582.2
$ws.Cells.Item(31, 12).Value = 7171.143
$ws.Cells.Item(31, 13).Value = -3227287.2
$ws.Cells.Item(31, 14).Value = -7761.143
# Row 34
$ws.Cells.Item(34, 8).Value = 2634348.5
$ws.Cells.Item(34, 9).Value = 3227582.2
$ws.Cells.Item(34, 10).Value = 7171.143
$ws.Cells.Item(34, 11).Value = 3227582.2
$ws.Cells.Item(34, 12).Value = 7171.143
$ws.Cells.Item(34, 13).Value = -3227380.2
$ws.Cells.Item(34, 14).Value = -7575.143
# Row 86
$ws.Cells.Item(86, 8).Value = 5955.6924
$ws.Cells.Item(86, 9).Value = 4053
$ws.Cells.Item(86, 10).Value = 9000
$ws.Cells.Item(86, 11).Value = 4053
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = -2930
$ws.Cells.Item(86, 14).Value = -11246
# Row 89
$ws.Cells.Item(89, 8).Value = 5955.6924
$ws.Cells.Item(89, 9).Value = 4053
$ws.Cells.Item(89, 10).Value = 9000
$ws.Cells.Item(89, 11).Value = 20265
$ws.Cells.Item(89, 12).Value = 45000
$ws.Cells.Item(89, 13).Value = -14649
$ws.Cells.Item(89, 14).Value = -56232
# Row 99
$ws.Cells.Item(99, 8).Value = 3052
$ws.Cells.Item(99, 9).Value = 1437.3334
$ws.Cells.Item(99, 10).Value = 4666.6665
$ws.Cells.Item(99, 11).Value = 1437.3334
$ws.Cells.Item(99, 12).Value = 4666.6665
$ws.Cells.Item(99, 13).Value = 60.66660000000002
$ws.Cells.Item(99, 14).Value = -7662.6665
# Row 122
$ws.Cells.Item(122, 8).Value = 4012.6
$ws.Cells.Item(122, 9).Value = 3795.8
$ws.Cells.Item(122, 10).Value = 4229.4
$ws.Cells.Item(122, 11).Value = 11387.4
$ws.Cells.Item(122, 12).Value = 12688.2
$ws.Cells.Item(122, 13).Value = -8937.400000000001
$ws.Cells.Item(122, 14).Value = -17588.2
# Row 126
$ws.Cells.Item(126, 8).Value = 3052
$ws.Cells.Item(126, 9).Value = 1437.3334
$ws.Cells.Item(126, 10).Value = 4666.6665
$ws.Cells.Item(126, 11).Value = 4312.0002
$ws.Cells.Item(126, 12).Value = 13999.9995
$ws.Cells.Item(126, 13).Value = -1842.0002
$ws.Cells.Item(126, 14).Value = -18939.9995
# Row 132
$ws.Cells.Item(132, 8).Value = 3501.7896
$ws.Cells.Item(132, 9).Value = 2688.2666
$ws.Cells.Item(132, 10).Value = 4032.348
$ws.Cells.Item(132, 11).Value = 8064.7998
$ws.Cells.Item(132, 12).Value = 12097.044
$ws.Cells.Item(132, 13).Value = -5534.7998
$ws.Cells.Item(132, 14).Value = -17157.044

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 889.8591300000001
$ws.Cells.Item(131, 10).Value = 1040.303
$ws.Cells.Item(131, 12).Value = 3120.909000000001
$ws.Cells.Item(131, 14).Value = -13200.909

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 2844.1843
$ws.Cells.Item(132, 9).Value = 2277.318
$ws.Cells.Item(132, 10).Value = 3623.625
$ws.Cells.Item(132, 11).Value = 6831.954000000001
$ws.Cells.Item(132, 12).Value = 10870.875
$ws.Cells.Item(132, 13).Value = -4301.954000000001
$ws.Cells.Item(132, 14).Value = -15930.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 2865
$ws.Cells.Item(7, 9).Value = 1950
$ws.Cells.Item(7, 11).Value = 1950
$ws.Cells.Item(7, 13).Value = -1838
# Row 61
$ws.Cells.Item(61, 8).Value = 2290.818
$ws.Cells.Item(61, 9).Value = 1519.9
$ws.Cells.Item(61, 11).Value = 1519.9
$ws.Cells.Item(61, 13).Value = -1317.9
# Row 100
$ws.Cells.Item(100, 8).Value = 2264.6155
$ws.Cells.Item(100, 9).Value = 1780
$ws.Cells.Item(100, 10).Value = 2567.5
$ws.Cells.Item(100, 11).Value = 1780
$ws.Cells.Item(100, 12).Value = 2567.5
$ws.Cells.Item(100, 13).Value = -1239
$ws.Cells.Item(100, 14).Value = -3649.5
# Row 113
$ws.Cells.Item(113, 8).Value = 2290.818
$ws.Cells.Item(113, 9).Value = 1519.9
$ws.Cells.Item(113, 11).Value = 1519.9
$ws.Cells.Item(113, 13).Value = 650.0999999999999
# Row 126
$ws.Cells.Item(126, 8).Value = 2865
$ws.Cells.Item(126, 9).Value = 1950
$ws.Cells.Item(126, 11).Value = 5850
$ws.Cells.Item(126, 13).Value = -3380

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 1395.1666
$ws.Cells.Item(113, 9).Value = 474.2
$ws.Cells.Item(113, 10).Value = 6000
$ws.Cells.Item(113, 11).Value = 1422.6
$ws.Cells.Item(113, 12).Value = 18000
$ws.Cells.Item(113, 13).Value = 747.4000000000001
$ws.Cells.Item(113, 14).Value = -22340

Write-Output "Applied 175 cell updates across 8 sheets."